$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 176, pushing the existing rows 176-182 down to 179-185.
$ws.Range("A176:A178").EntireRow.Insert()

# Fill the 3 newly inserted rows (176-178) with the new "Doctor Davis" entries.
# Columns A,B,C,E,F,G,H,I,J,R are identical across all rows in this block.

# Row 176
$ws.Range("A176").Value2 = 11
$ws.Range("B176").Value2 = "Vega Monumental Concepción"
$ws.Range("C176").Value2 = "Bíobío"
$ws.Range("D176").Value2 = 44615
$ws.Range("E176").Value2 = 8
$ws.Range("F176").Value2 = "Fruta"
$ws.Range("G176").Value2 = 100103
$ws.Range("H176").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I176").Value2 = 100103004
$ws.Range("J176").Value2 = "Durazno"
$ws.Range("K176").Value2 = "Doctor Davis"
$ws.Range("L176").Value2 = "Especial"
$ws.Range("M176").Value2 = 50
$ws.Range("N176").Value2 = 13000
$ws.Range("O176").Value2 = 13000
$ws.Range("P176").Value2 = 13000
$ws.Range("Q176").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R176").Value2 = "Región de O'Higgins"
$ws.Range("S176").Value2 = 867
$ws.Range("T176").Value2 = 15

# Row 177
$ws.Range("A177").Value2 = 11
$ws.Range("B177").Value2 = "Vega Monumental Concepción"
$ws.Range("C177").Value2 = "Bíobío"
$ws.Range("D177").Value2 = 44615
$ws.Range("E177").Value2 = 8
$ws.Range("F177").Value2 = "Fruta"
$ws.Range("G177").Value2 = 100103
$ws.Range("H177").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I177").Value2 = 100103004
$ws.Range("J177").Value2 = "Durazno"
$ws.Range("K177").Value2 = "Doctor Davis"
$ws.Range("L177").Value2 = "Primera"
$ws.Range("M177").Value2 = 100
$ws.Range("N177").Value2 = 11000
$ws.Range("O177").Value2 = 11000
$ws.Range("P177").Value2 = 11000
$ws.Range("Q177").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R177").Value2 = "Región de O'Higgins"
$ws.Range("S177").Value2 = 733
$ws.Range("T177").Value2 = 15

# Row 178
$ws.Range("A178").Value2 = 11
$ws.Range("B178").Value2 = "Vega Monumental Concepción"
$ws.Range("C178").Value2 = "Bíobío"
$ws.Range("D178").Value2 = 44615
$ws.Range("E178").Value2 = 8
$ws.Range("F178").Value2 = "Fruta"
$ws.Range("G178").Value2 = 100103
$ws.Range("H178").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I178").Value2 = 100103004
$ws.Range("J178").Value2 = "Durazno"
$ws.Range("K178").Value2 = "Doctor Davis"
$ws.Range("L178").Value2 = "Segunda"
$ws.Range("M178").Value2 = 100
$ws.Range("N178").Value2 = 9000
$ws.Range("O178").Value2 = 9000
$ws.Range("P178").Value2 = 9000
$ws.Range("Q178").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R178").Value2 = "Región de O'Higgins"
$ws.Range("S178").Value2 = 600
$ws.Range("T178").Value2 = 15
